$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inspection Log")
$lo = $ws.ListObjects.Item(1)
for ($i=0; $i -lt 4; $i++) {
  $lo.ListRows.Add() | Out-Null
}
$rng = $ws.Range("A162:E165")
$white = 16777215
# xlInsideVertical=11, xlInsideHorizontal=12, xlEdgeLeft=7,xlEdgeTop=8,xlEdgeBottom=9,xlEdgeRight=10
foreach ($idx in 7,8,9,10,11,12) {
  $rng.Borders.Item($idx).LineStyle = 1
  $rng.Borders.Item($idx).Weight = 2
  $rng.Borders.Item($idx).Color = $white
}
Write-Output "ok"
